$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '95.282.67'
$ws.Range("E2").Value = '  +2.11%  '

# Row 3
$ws.Range("D3").Value = '3.583.61'
$ws.Range("E3").Value = '  +4.72%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.10'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.81%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '652.41'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +5.08%  '

# Row 7
$ws.Range("E7").Value = '  +6.44%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.405'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +3.02%  '

# Row 10
$ws.Range("E10").Value = '  +3.86%  '

# Row 11
$ws.Range("D11").Value = '3.581.01'
$ws.Range("E11").Value = '  +4.70%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.95'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.33%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.200'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.12%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.29'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.68%  '

# Row 15
$ws.Range("D15").Value = '4.265.63'
$ws.Range("E15").Value = '  +5.17%  '

# Row 16
$ws.Range("D16").Value = '95.036.49'
$ws.Range("E16").Value = '  +2.06%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000255'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +3.50%  '

# Row 18
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.592.45'
$ws.Range("E18").Value = '  +5.12%  '

# Row 19
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.93'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.12%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.57'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +7.64%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.97'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.94%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.58'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +7.25%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.484'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +8.33%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '511.36'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.47%  '

# Row 25
$ws.Range("E25").Value = '  +5.32%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.62'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.00%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '96.56'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.74%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.68'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +5.86%  '

# Row 29
$ws.Range("D29").Value = '3.778.85'
$ws.Range("E29").Value = '  +4.94%  '

# Row 30
$ws.Range("E30").Value = '  +17.18%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.27'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.44%  '

# Row 32
$ws.Range("E32").Value = '  -0.13%  '

# Row 33
$ws.Range("E33").Value = '  +1.45%  '

# Row 34
$ws.Range("E34").Value = '  -0.61%  '

# Row 35
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '31.91'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +6.33%  '

# Row 36
$ws.Range("B36").Value = 'Cronos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.176'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.60%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.557'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.07%  '

# Row 38
$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '576.87'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +4.94%  '

# Row 39
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.22'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +10.00%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.49'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +5.79%  '

# Row 41
$ws.Range("E41").Value = '  -0.03%  '

# Row 42
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.926'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.46%  '

# Row 43
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.150'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.55%  '

# Row 44
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.73'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.61%  '

# Row 45
$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '23.74'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.27%  '

# Row 46
$ws.Range("B46").Value = 'ImmutableX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.71'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.75%  '

# Row 47
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.25'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +6.26%  '

# Row 48
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0417'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.82%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.01'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +30.62%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.97'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.50%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.47'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -6.41%  '
